# Applies the env_site_final.xlsx edit:
#   - Column K (ASP, aspect in degrees): mirror the value about 90 deg,
#     i.e. new = 180 - old, for every data row EXCEPT the "no data"
#     sentinel rows where the stored aspect is exactly 0.1.
#   - Column Q (Moose_new): decrement by 1 for every row belonging to the
#     "GM" park, except rows already at the floor value of 1.
#
# Columns by letter -> index: F=6 (park), K=11 (ASP), Q=17 (Moose_new)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $kCell = $ws.Cells.Item($r, 11)
    $kVal = $kCell.Value()
    if ($kVal -ne $null -and $kVal -ne "" -and $kVal -ne 0.1) {
        $kCell.Value = 180 - $kVal
    }

    $fCell = $ws.Cells.Item($r, 6)
    $fVal = $fCell.Value()
    if ($fVal -eq "GM") {
        $qCell = $ws.Cells.Item($r, 17)
        $qVal = $qCell.Value()
        if ($qVal -ne $null -and $qVal -ne "" -and $qVal -gt 1) {
            $qCell.Value = $qVal - 1
        }
    }
}
